$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A171").Value = "IMX-USD"
$ws.Range("A172").Value = "TAO-USD"
$ws.Range("A173").Value = "GRT-USD"
